$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for years 2004-2009 (rows 2-7), shifting 2010-2020 data up to rows 2-12
$ws.Rows("2:7").Delete()

# Add new row 13 for year 2021, copying the style of column A from the row above
$ws.Cells.Item(12, 1).Copy()
$ws.Cells.Item(13, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 185.03
$ws.Cells.Item(13, 3).Value = 673.344245
$ws.Cells.Item(13, 4).Value = 1059901.18
$ws.Cells.Item(13, 5).Value = 31737.67
$ws.Cells.Item(13, 6).Value = 170.436268
$ws.Cells.Item(13, 7).Value = 375.378253
$ws.Cells.Item(13, 8).Value = 55580.86
